$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row text fix-up ---------------------------------------------
# Column C was mislabeled "username" (it has always held the registered
# e-mail address, see the hyperlink fix below) -> relabel to "email".
# FirstName/LastName -> firstName/lastName (camelCase to match the other
# automation data-provider keys). "password" (column D) is left as-is.
$ws.Range("A1").Value = "firstName"
$ws.Range("B1").Value = "lastName"
$ws.Range("C1").Value = "email"

# --- Header styling: shrink the oversized header font -------------------
# Header stays bold/filled/centered/bordered, just drops from 14pt to the
# normal 11pt used by the rest of the workbook.
$ws.Range("A1:D1").Font.Size = 11
# Row 1 had an explicit 18pt height to fit the old 14pt font; let Excel
# recompute the (now default) row height instead of leaving it pinned.
$ws.Rows(1).AutoFit() | Out-Null

# --- Hyperlinks: the two mailto links were attached to the wrong cells --
# C2 displays the e-mail address and should link to it; D2 displays the
# password text and should link to the "Test@1234" mailto. Re-create both
# links with their targets swapped back to the correct display cell.
$ws.Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:Auto.Reg1rew@yopmail.com", "", "", "Test@1234") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Test@1234", "", "", "Auto.Reg1rew@yopmail.com") | Out-Null

# --- Selection: author left column C selected when the file was saved ---
$ws.Columns("C").Select() | Out-Null
